$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 267 (shifts old rows 267-308 down to 268-309),
# matching the weekly data-refresh pattern: the newest week's record is
# prepended to the historical list right after the previous week's entry.
$ws.Rows.Item(267).Insert()

$ws.Range("A267").Value = 11
$ws.Range("B267").Value = "Vega Monumental Concepción"
$ws.Range("C267").Value = "Bíobío"
$ws.Range("D267").Value = 45127
$ws.Range("E267").Value = 8
$ws.Range("F267").Value = 100112003
$ws.Range("G267").Value = "Ajo"
$ws.Range("H267").Value = "Chino"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 220
$ws.Range("K267").Value = 15000
$ws.Range("L267").Value = 16000
$ws.Range("M267").Value = 15545
$ws.Range("N267").Value = "$/caja 10 kilos"
$ws.Range("O267").Value = "China"
$ws.Range("P267").Value = 1554
$ws.Range("Q267").Value = 10
$ws.Range("R267").Value = "Hortaliza"
